# Non-Oncology latest template changes
#
# The "scenario6" error-message list (column F, rows 18-29) drops the
# "Duplicate column found for QA-3..." entry that used to live on row 25.
# Deleting that entire row shifts rows 26-29 up to 25-28, which also
# removes the now-unused shared string from the table automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(25).Delete()

# Column E (Files_to_upload) was widened / best-fit to accommodate the
# long template file paths.
$ws.Columns("E").ColumnWidth = 113.5

# The saved view now has the window scrolled down with row 25 selected
# (full-row selection, A25:XFD25).
$ws.Range("A25:XFD25").Select()
